$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("WebViewMeasurement")

# --- Window / workbook-level cosmetic props (best effort) ---
$win = $excel.ActiveWindow
$win.Width = 16200

# --- Clear the two long commentary text cells (keep style, drop content) ---
$ws1.Range("A1").Value = $null
$ws1.Range("A23").Value = $null

# --- Fill in first table (rows 10-19) with new measured results ---
$ws1.Range("B10").Value = 8.4
$ws1.Range("C10").Value = 5.28
$ws1.Range("E10").Value = 8.85
$ws1.Range("F10").Value = 8.86

$ws1.Range("B11").Value = 4.58
$ws1.Range("C11").Value = 4.98
$ws1.Range("E11").Value = 8.41
$ws1.Range("F11").Value = 8.52

$ws1.Range("B12").Value = 5.35
$ws1.Range("C12").Value = 4.96
$ws1.Range("E12").Value = 8.59
$ws1.Range("F12").Value = 9.09

$ws1.Range("B13").Value = 5.35
$ws1.Range("C13").Value = 5.23

$ws1.Range("B14").Value = 5.34
$ws1.Range("C14").Value = 5.21

$ws1.Range("B15").Value = 8.25
$ws1.Range("C15").Value = 4.16

$ws1.Range("B16").Value = 8.28
$ws1.Range("C16").Value = 5.09

$ws1.Range("B17").Value = 9.61
$ws1.Range("C17").Value = 4.67

$ws1.Range("B18").Value = 5.22
$ws1.Range("C18").Value = 5.07

$ws1.Range("B19").Value = 6.62
$ws1.Range("C19").Value = 5.38

# --- Clear out second table (rows 33-42), data moved/removed ---
$ws1.Range("B33:C33").Value = $null
$ws1.Range("E33:F33").Value = $null
$ws1.Range("B34:C34").Value = $null
$ws1.Range("E34:F34").Value = $null
$ws1.Range("B35:C35").Value = $null
$ws1.Range("E35:F35").Value = $null
$ws1.Range("B36:C36").Value = $null
$ws1.Range("E36:F36").Value = $null
$ws1.Range("B37:C37").Value = $null
$ws1.Range("E37:F37").Value = $null
$ws1.Range("B38:C38").Value = $null
$ws1.Range("E38:F38").Value = $null
$ws1.Range("B39:C39").Value = $null
$ws1.Range("E39:F39").Value = $null
$ws1.Range("B40:C40").Value = $null
$ws1.Range("B41:C41").Value = $null
$ws1.Range("B42:C42").Value = $null

# --- Selection moved to F13 ---
$ws1.Range("F13").Select()

Write-Output "edit complete"
